# Generate Report for Handback
# Refresh the handoff/handback timestamps and priority for the
# c95d6e5e-... and df100312-... entries across the Overview, zh-cn and
# de-de sheets, as produced by a re-run of the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for c95d6e5e-...md (row 4)
# and df100312-...md (row 5) both move from 02:16:02 to 02:16:54.
$wsOverview.Range("G4").Value = "2016-08-28 02:16:54"
$wsOverview.Range("G5").Value = "2016-08-28 02:16:54"

# zh-cn table: Priority flips from "ht" (human translated) to "mt"
# (machine translated) for both rows.
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn table: Correspond Handoff / Handback datetimes move forward.
$wsZhCn.Range("H4").Value = "2016-08-28 02:16:49"
$wsZhCn.Range("H5").Value = "2016-08-28 02:16:49"
$wsZhCn.Range("K4").Value = "2016-08-28 02:17:12"
$wsZhCn.Range("K5").Value = "2016-08-28 02:17:12"

# de-de table: Priority flips from "ht" to "mt" for both rows.
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de table: Correspond Handoff datetime moves forward (shared with
# the Overview "Latest HO Xliff Generate Date" string).
$wsDeDe.Range("H4").Value = "2016-08-28 02:16:54"
$wsDeDe.Range("H5").Value = "2016-08-28 02:16:54"

# de-de table: Correspond Handback datetime moves forward.
$wsDeDe.Range("K4").Value = "2016-08-28 02:17:19"
$wsDeDe.Range("K5").Value = "2016-08-28 02:17:19"
